$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "Owner" column (C) values for the relevant rows, and update the
# "Implementation/Design Notes" text for the Strategy pattern row (row 2).
$ws.Range("C2").Value = "Joe"
$ws.Range("E2").Value = "An interface ProgressibleUnit that has a method called progress() to handle everything each unit needs to do with each game iterration."
$ws.Range("C3").Value = "Vince"
$ws.Range("C4").Value = "Zach"
$ws.Range("C6").Value = "Vince"
$ws.Range("C7").Value = "Rocky"
$ws.Range("C8").Value = "Zach"

# Match the style used by the rest of the row (column B) for the new C data cells:
# horizontally and vertically centered, like the other owner/pattern cells.
$newCells = "C2", "C3", "C4", "C6", "C7", "C8"
foreach ($addr in $newCells) {
    $ws.Range($addr).HorizontalAlignment = $excel.Constants.xlCenter
    $ws.Range($addr).VerticalAlignment = $excel.Constants.xlCenter
}

# Row 2 height grew to accommodate the longer owner/notes text
$ws.Rows("2").RowHeight = 34

# Update the selected cell on the frozen pane to C8
$ws.Range("C8").Select()
